# employeesDB.xlsx update
#  - Remove the "luis / 22da" record (old row 3), shifting the rows below it up.
#  - Append a new employee record "Amarilis / marca" with Salary/Department 12321.0.
#
# Resulting data:
#   Name     | Lastname  | Salary  | Department
#   luis     | javier    | 2341.0  | 2341.0
#   javier   | tatis     | 12322.0 | 12322.0
#   maria    | maldonado | 12212.0 | 12212.0
#   Amarilis | marca     | 12321.0 | 12321.0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the obsolete "luis / 22da" row; rows below shift up automatically.
$ws.Rows.Item(3).Delete()

# Copy the formatting of the (now last) data row down onto the new row
# so the appended record keeps the same style used by the other records.
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)

# Populate the new trailing record.
$ws.Range("A5").Value = "Amarilis"
$ws.Range("B5").Value = "marca"
$ws.Range("C5").Value = 12321.0
$ws.Range("D5").Value = 12321.0
